$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54 values
$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -0.7200474048664085
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -2.055952042396259

# Copy the formatting (date number format, bold font, border, alignment)
# from the row above (A53) so the new date cell A54 matches the existing
# column-A style used throughout the sheet.
$ws.Range("A53").Copy() | Out-Null
$ws.Range("A54").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
